$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(132, 8).Value = 1900.1923
$ws.Cells.Item(132, 9).Value = 1304.1
$ws.Cells.Item(132, 11).Value = 3912.3
$ws.Cells.Item(132, 13).Value = -1382.3

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(12, 8).Value = 300
$ws.Cells.Item(12, 10).Value = 0
$ws.Cells.Item(12, 12).Value = 0
$ws.Cells.Item(12, 14).ClearContents()
$ws.Cells.Item(31, 8).Value = 17078.445
$ws.Cells.Item(31, 9).Value = 7123.4287
$ws.Cells.Item(31, 11).Value = 7123.4287
$ws.Cells.Item(31, 13).Value = -6829.4287
$ws.Cells.Item(63, 8).Value = 2850
$ws.Cells.Item(63, 10).Value = 9000
$ws.Cells.Item(63, 12).Value = 9000
$ws.Cells.Item(63, 14).Value = -10372
$ws.Cells.Item(66, 8).Value = 2850
$ws.Cells.Item(66, 10).Value = 9000
$ws.Cells.Item(66, 12).Value = 45000
$ws.Cells.Item(66, 14).Value = -51864
$ws.Cells.Item(74, 8).Value = 31818.941
$ws.Cells.Item(74, 10).Value = 4555.5557
$ws.Cells.Item(74, 12).Value = 4555.5557
$ws.Cells.Item(74, 14).Value = -6303.5557
$ws.Cells.Item(77, 8).Value = 31818.941
$ws.Cells.Item(77, 10).Value = 4555.5557
$ws.Cells.Item(77, 12).Value = 22777.7785
$ws.Cells.Item(77, 14).Value = -31513.7785
$ws.Cells.Item(122, 8).Value = 3222.394
$ws.Cells.Item(122, 9).Value = 1315.579
$ws.Cells.Item(122, 11).Value = 3946.737
$ws.Cells.Item(122, 13).Value = -1496.737
$ws.Cells.Item(134, 8).Value = 30000
$ws.Cells.Item(134, 10).Value = 30000
$ws.Cells.Item(134, 12).Value = 30000
$ws.Cells.Item(134, 14).Value = -40140

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(92, 8).Value = 0
$ws.Cells.Item(92, 10).Value = 0
$ws.Cells.Item(92, 12).ClearContents()
$ws.Cells.Item(92, 14).Value = 0
$ws.Cells.Item(134, 8).Value = 6027.825
$ws.Cells.Item(134, 9).Value = 2684.3333
$ws.Cells.Item(134, 11).Value = 8052.999899999999
$ws.Cells.Item(134, 13).Value = -5517.999899999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 6890.5654
$ws.Cells.Item(31, 9).Value = 1778.25
$ws.Cells.Item(31, 10).Value = 10823.115
$ws.Cells.Item(31, 11).Value = 1778.25
$ws.Cells.Item(31, 12).Value = 10823.115
$ws.Cells.Item(31, 13).Value = -1483.25
$ws.Cells.Item(31, 14).Value = -11413.115
$ws.Cells.Item(34, 8).Value = 6890.5654
$ws.Cells.Item(34, 9).Value = 1778.25
$ws.Cells.Item(34, 10).Value = 10823.115
$ws.Cells.Item(34, 11).Value = 1778.25
$ws.Cells.Item(34, 12).Value = 10823.115
$ws.Cells.Item(34, 13).Value = -1576.25
$ws.Cells.Item(34, 14).Value = -11227.115

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(17, 8).Value = 1294.75
$ws.Cells.Item(17, 9).Value = 342
$ws.Cells.Item(17, 10).Value = 2882.6667
$ws.Cells.Item(17, 11).Value = 1026
$ws.Cells.Item(17, 12).Value = 8648.000100000001
$ws.Cells.Item(17, 13).Value = -857
$ws.Cells.Item(17, 14).Value = -8986.000100000001
$ws.Cells.Item(131, 8).Value = 3223.2104
$ws.Cells.Item(131, 9).Value = 1626.25
$ws.Cells.Item(131, 11).Value = 4878.75
$ws.Cells.Item(131, 13).Value = 161.25
$ws.Cells.Item(136, 8).Value = 2409.2856
$ws.Cells.Item(136, 9).Value = 2409.2856
$ws.Cells.Item(136, 11).Value = 7227.8568
$ws.Cells.Item(136, 13).Value = -2127.8568
$ws.Cells.Item(137, 8).Value = 156150.47
$ws.Cells.Item(137, 9).Value = 112967.555
$ws.Cells.Item(137, 11).Value = 338902.665
$ws.Cells.Item(137, 13).Value = -333802.665
$ws.Cells.Item(139, 8).Value = 68134.75
$ws.Cells.Item(139, 9).Value = 103015.7
$ws.Cells.Item(139, 11).Value = 309047.1
$ws.Cells.Item(139, 13).Value = -303907.1
$ws.Cells.Item(140, 8).Value = 88270.61
$ws.Cells.Item(140, 9).Value = 88270.61
$ws.Cells.Item(140, 11).Value = 264811.83
$ws.Cells.Item(140, 13).Value = -259631.83

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(12, 8).Value = 200
$ws.Cells.Item(12, 10).Value = 200
$ws.Cells.Item(12, 12).Value = 200
$ws.Cells.Item(12, 14).Value = -480
$ws.Cells.Item(14, 8).Value = 0
$ws.Cells.Item(14, 9).Value = 0
$ws.Cells.Item(14, 11).Value = 0
$ws.Cells.Item(14, 13).ClearContents()
$ws.Cells.Item(43, 8).Value = 1891.8334
$ws.Cells.Item(43, 9).Value = 1891.8334
$ws.Cells.Item(43, 11).Value = 1891.8334
$ws.Cells.Item(43, 13).Value = -1740.8334
$ws.Cells.Item(113, 8).Value = 6980.7427
$ws.Cells.Item(113, 9).Value = 4161.0557
$ws.Cells.Item(113, 10).Value = 9966.294
$ws.Cells.Item(113, 11).Value = 4161.0557
$ws.Cells.Item(113, 12).Value = 9966.294
$ws.Cells.Item(113, 13).Value = -1991.0557
$ws.Cells.Item(113, 14).Value = -14306.294
$ws.Cells.Item(126, 8).Value = 4699.1763
$ws.Cells.Item(126, 9).Value = 2662.25
$ws.Cells.Item(126, 10).Value = 6509.778
$ws.Cells.Item(126, 11).Value = 7986.75
$ws.Cells.Item(126, 12).Value = 19529.334
$ws.Cells.Item(126, 13).Value = -5516.75
$ws.Cells.Item(126, 14).Value = -24469.334
$ws.Cells.Item(132, 8).Value = 4654.793
$ws.Cells.Item(132, 9).Value = 2652.1738
$ws.Cells.Item(132, 10).Value = 12331.5
$ws.Cells.Item(132, 11).Value = 7956.5214
$ws.Cells.Item(132, 12).Value = 36994.5
$ws.Cells.Item(132, 13).Value = -5426.5214
$ws.Cells.Item(132, 14).Value = -42054.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(22, 8).Value = 1559.4615
$ws.Cells.Item(22, 9).Value = 707.44446
$ws.Cells.Item(22, 10).Value = 3476.5
$ws.Cells.Item(22, 11).Value = 707.44446
$ws.Cells.Item(22, 12).Value = 3476.5
$ws.Cells.Item(22, 13).Value = -412.44446
$ws.Cells.Item(22, 14).Value = -4066.5
$ws.Cells.Item(27, 8).Value = 1559.4615
$ws.Cells.Item(27, 9).Value = 707.44446
$ws.Cells.Item(27, 10).Value = 3476.5
$ws.Cells.Item(27, 11).Value = 707.44446
$ws.Cells.Item(27, 12).Value = 3476.5
$ws.Cells.Item(27, 13).Value = -600.44446
$ws.Cells.Item(27, 14).Value = -3690.5
$ws.Cells.Item(46, 8).Value = 2851.65
$ws.Cells.Item(46, 9).Value = 2266.6428
$ws.Cells.Item(46, 10).Value = 4216.6665
$ws.Cells.Item(46, 11).Value = 2266.6428
$ws.Cells.Item(46, 12).Value = 4216.6665
$ws.Cells.Item(46, 13).Value = -2078.6428
$ws.Cells.Item(46, 14).Value = -4592.6665
$ws.Cells.Item(55, 8).Value = 362.4516
$ws.Cells.Item(55, 9).Value = 109.066666
$ws.Cells.Item(55, 10).Value = 600
$ws.Cells.Item(55, 11).Value = 109.066666
$ws.Cells.Item(55, 12).Value = 600
$ws.Cells.Item(55, 13).Value = 63.933334
$ws.Cells.Item(55, 14).Value = -946
$ws.Cells.Item(68, 8).Value = 6065.846
$ws.Cells.Item(68, 9).Value = 3891.4
$ws.Cells.Item(68, 10).Value = 7424.875
$ws.Cells.Item(68, 11).Value = 3891.4
$ws.Cells.Item(68, 12).Value = 7424.875
$ws.Cells.Item(68, 13).Value = -3142.4
$ws.Cells.Item(68, 14).Value = -8922.875
$ws.Cells.Item(71, 8).Value = 6065.846
$ws.Cells.Item(71, 9).Value = 3891.4
$ws.Cells.Item(71, 10).Value = 7424.875
$ws.Cells.Item(71, 11).Value = 19457
$ws.Cells.Item(71, 12).Value = 37124.375
$ws.Cells.Item(71, 13).Value = -15713
$ws.Cells.Item(71, 14).Value = -44612.375
$ws.Cells.Item(136, 8).Value = 13707.632
$ws.Cells.Item(136, 9).Value = 3052.5293
$ws.Cells.Item(136, 10).Value = 22333.191
$ws.Cells.Item(136, 11).Value = 9157.5879
$ws.Cells.Item(136, 12).Value = 66999.573
$ws.Cells.Item(136, 13).Value = -6607.5879
$ws.Cells.Item(136, 14).Value = -72099.573

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(14, 8).Value = 250253250
$ws.Cells.Item(14, 9).Value = 250253250
$ws.Cells.Item(14, 10).Value = 0
$ws.Cells.Item(14, 11).Value = 250253250
$ws.Cells.Item(14, 12).Value = 0
$ws.Cells.Item(14, 13).ClearContents()
$ws.Cells.Item(14, 14).Value = -250253082
$ws.Cells.Item(62, 8).Value = 166671460
$ws.Cells.Item(62, 10).Value = 166668670
$ws.Cells.Item(62, 12).Value = 166668670
$ws.Cells.Item(62, 14).Value = -166669918
$ws.Cells.Item(65, 8).Value = 166671460
$ws.Cells.Item(65, 10).Value = 166668670
$ws.Cells.Item(65, 12).Value = 833343350
$ws.Cells.Item(65, 14).Value = -833349590
$ws.Cells.Item(100, 8).Value = 886
$ws.Cells.Item(100, 10).Value = 1024.25
$ws.Cells.Item(100, 12).Value = 2048.5
$ws.Cells.Item(100, 14).Value = -3130.5
$ws.Cells.Item(132, 8).Value = 5375.1787
$ws.Cells.Item(132, 9).Value = 7443.7144
$ws.Cells.Item(132, 10).Value = 3306.6428
$ws.Cells.Item(132, 11).Value = 22331.1432
$ws.Cells.Item(132, 12).Value = 9919.928400000001
$ws.Cells.Item(132, 13).Value = -19801.1432
$ws.Cells.Item(132, 14).Value = -14979.9284
$ws.Cells.Item(136, 8).Value = 390085.47
$ws.Cells.Item(136, 10).Value = 632725.2
$ws.Cells.Item(136, 12).Value = 1898175.6
$ws.Cells.Item(136, 14).Value = -1903275.6
